# Week 17 data logging + tiebreak fix for "2021 Team Data.xlsx"
# Appends the week's per-play logs onto the running season strings and
# updates the aggregate total cells on the OFF / DEF / ST / TURNS / PEN sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# YDS sheet: append this week's play-by-play yardage logs
# ---------------------------------------------------------------------------
$wsYDS = $wb.Worksheets.Item("YDS")

$wsYDS.Range("B2").Value = $wsYDS.Range("B2").Text + " 2 4 3 0 4 4 13 9 -4 8 6 4 1 11 -3 7 5 2 13 0 6"
$wsYDS.Range("B3").Value = $wsYDS.Range("B3").Text + " 8 15 0 9 10 61 15 24 4 1 8 15 17 11"
$wsYDS.Range("C2").Value = $wsYDS.Range("C2").Text + " 6 0 1 1 3 4 3 20 1 16 6 6 4 12 3 5 20 9 5 -2 14 4 7 3 6 11 5 5 4 2 2 12 4 4 6 5 4 2 3 8"
$wsYDS.Range("C3").Value = $wsYDS.Range("C3").Text + " 13 15 11 15 5 6 10 13 9 15 8 2"

# ---------------------------------------------------------------------------
# OFF sheet: updated season totals after Week 17
# ---------------------------------------------------------------------------
$wsOFF = $wb.Worksheets.Item("OFF")

$wsOFF.Range("C2").Value = 208
$wsOFF.Range("D2").Value = 6
$wsOFF.Range("F2").Value = 35
$wsOFF.Range("G2").Value = 42
$wsOFF.Range("N2").Value = 18
$wsOFF.Range("O2").Value = 21

$wsOFF.Range("C3").Value = 139
$wsOFF.Range("E3").Value = 29
$wsOFF.Range("F3").Value = 131
$wsOFF.Range("G3").Value = 35
$wsOFF.Range("H3").Value = 26
$wsOFF.Range("I3").Value = 79
$wsOFF.Range("J3").Value = 50
$wsOFF.Range("L3").Value = 287
$wsOFF.Range("M3").Value = 189
$wsOFF.Range("Q3").Value = 522

# ---------------------------------------------------------------------------
# DEF sheet: updated season totals after Week 17
# ---------------------------------------------------------------------------
$wsDEF = $wb.Worksheets.Item("DEF")

$wsDEF.Range("B2").Value = 14
$wsDEF.Range("C2").Value = 199
$wsDEF.Range("D2").Value = 10
$wsDEF.Range("E2").Value = 12
$wsDEF.Range("F2").Value = 64
$wsDEF.Range("G2").Value = 76
$wsDEF.Range("J2").Value = 40
$wsDEF.Range("N2").Value = 10
$wsDEF.Range("O2").Value = 21
$wsDEF.Range("P2").Value = 11

$wsDEF.Range("C3").Value = 181
$wsDEF.Range("D3").Value = 10
$wsDEF.Range("E3").Value = 31
$wsDEF.Range("F3").Value = 103
$wsDEF.Range("G3").Value = 31
$wsDEF.Range("I3").Value = 52
$wsDEF.Range("J3").Value = 54
$wsDEF.Range("L3").Value = 312
$wsDEF.Range("M3").Value = 204
$wsDEF.Range("Q3").Value = 579

# ---------------------------------------------------------------------------
# ST sheet: updated season totals + appended per-kick logs after Week 17
# ---------------------------------------------------------------------------
$wsST = $wb.Worksheets.Item("ST")

$wsST.Range("B2").Value = 62
$wsST.Range("D2").Value = 63
$wsST.Range("F2").Value = 89
$wsST.Range("G2").Value = 85
$wsST.Range("J2").Value = 59
$wsST.Range("K2").Value = 57

$wsST.Range("B4").Value = $wsST.Range("B4").Text + " 65 66 61"
$wsST.Range("B5").Value = $wsST.Range("B5").Text + " 23 26 25"
$wsST.Range("B6").Value = $wsST.Range("B6").Text + " 23"
$wsST.Range("D3").Value = $wsST.Range("D3").Text + " 34 34 40 37 48"
$wsST.Range("D4").Value = $wsST.Range("D4").Text + " 0 0 0 0 19"

# ---------------------------------------------------------------------------
# TURNS sheet: updated season totals after Week 17
# ---------------------------------------------------------------------------
$wsTURNS = $wb.Worksheets.Item("TURNS")

$wsTURNS.Range("C3").Value = 8
$wsTURNS.Range("D3").Value = 12

# ---------------------------------------------------------------------------
# PEN sheet: updated season totals after Week 17
# ---------------------------------------------------------------------------
$wsPEN = $wb.Worksheets.Item("PEN")

$wsPEN.Range("D4").Value = 11
